# "Version 1.1.4 - Staging"
#
# 1) Bump the cached "datetimeFigureOut" footer-date field (6/4/25 -> 6/5/25)
#    everywhere it is cached: the slide master and every slide layout.
# 2) Slide 3 ("Review Testing"): merge the 3 runs of the "Create Course
#    Scrolling Issue ..." bullet into a single run.
# 3) Slide 8 ("Next STEPS"): bump "Update to Version 1.0.9" -> "1.1.2".
#
# NOTE: this host's PowerShell functions do not get their own variable
# scope (loop counters leak into/clobber the caller's), so every
# function/loop below uses its own uniquely-named counter variable to
# avoid accidentally resetting an outer loop's counter.

function Find-ShapeByNameContains($shapes, $substr) {
    for ($fsbnc_i = 1; $fsbnc_i -le $shapes.Count; $fsbnc_i++) {
        $fsbnc_sh = $shapes.Item($fsbnc_i)
        if ($fsbnc_sh.Name -like "*$substr*") {
            return $fsbnc_sh
        }
    }
    return $null
}

# Force a real text-range mutation: this host's engine treats an
# assignment that nets out to the same display string as a no-op, which
# would leave multiple runs / the <a:fld> field intact instead of the
# single collapsed run PowerPoint itself produces when you retype text
# over it. Writing a throwaway value first guarantees the later write
# lands as one fresh run.
function Set-RangeText($srt_range, [string]$srt_text) {
    $srt_range.Text = "~~tmp~~"
    $srt_range.Text = $srt_text
}

$p = $ppt.ActivePresentation

# --- 1) Footer date field: slide master + all slide layouts ---

$masterDateShape = Find-ShapeByNameContains $p.SlideMaster.Shapes "Date"
if ($masterDateShape -ne $null) {
    Set-RangeText $masterDateShape.TextFrame.TextRange "6/5/25"
}

$layouts = $p.SlideMaster.CustomLayouts
for ($layoutIdx = 1; $layoutIdx -le $layouts.Count; $layoutIdx++) {
    $layout = $layouts.Item($layoutIdx)
    $layoutDateShape = Find-ShapeByNameContains $layout.Shapes "Date"
    if ($layoutDateShape -ne $null) {
        Set-RangeText $layoutDateShape.TextFrame.TextRange "6/5/25"
    }
}

# --- 2) Slide 3: merge the "Create Course Scrolling Issue" runs ---

$slide3 = $p.Slides.Item(3)
$contentShape3 = $slide3.Shapes.Item(2)
$contentRange3 = $contentShape3.TextFrame.TextRange
$contentParas3 = $contentRange3.Paragraphs()
for ($p3i = 1; $p3i -le $contentParas3.Count; $p3i++) {
    $para3 = $contentRange3.Paragraphs($p3i, 1)
    if ($para3.Text -like "Create Course Scrolling Issue*") {
        Set-RangeText $para3 "Create Course Scrolling Issue – Possibly Fixed in Version 1.1.3"
        break
    }
}

# --- 3) Slide 8: version bump 1.0.9 -> 1.1.2 ---

$slide8 = $p.Slides.Item(8)
$contentShape8 = $slide8.Shapes.Item(2)
$contentRange8 = $contentShape8.TextFrame.TextRange
$contentParas8 = $contentRange8.Paragraphs()
for ($p8i = 1; $p8i -le $contentParas8.Count; $p8i++) {
    $para8 = $contentRange8.Paragraphs($p8i, 1)
    if ($para8.Text -like "Update to Version*") {
        Set-RangeText $para8 "Update to Version 1.1.2"
        break
    }
}

Write-Host "Done."
